$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Seed the shared-string table in the same registration order as the
# source workbook: "5:37:42 PM" must become the first newly-added unique
# string, "5:32:51 PM" the second — so set F95's text before F94's.
$ws.Cells.Item(95, 6).Value = "Apr 25, 2024 5:37:42 PM"
$ws.Cells.Item(94, 6).Value = "Apr 25, 2024 5:32:51 PM"

# Add row 94 (A94 = 93)
$ws.Cells.Item(94, 1).Value = 93
$ws.Cells.Item(94, 2).Value = "patrick.fernandes"
$ws.Cells.Item(94, 3).Value = "?"
$ws.Cells.Item(94, 4).Value = "TOTVS Printer Document"
$ws.Cells.Item(94, 5).Value = "Impressão concluída"
$ws.Cells.Item(94, 7).Value = 11

# Add row 95 (A95 = 94)
$ws.Cells.Item(95, 1).Value = 94
$ws.Cells.Item(95, 2).Value = "patrick.fernandes"
$ws.Cells.Item(95, 3).Value = "?"
$ws.Cells.Item(95, 4).Value = "TOTVS Printer Document"
$ws.Cells.Item(95, 5).Value = "Impressão concluída"
$ws.Cells.Item(95, 7).Value = 1

# Update selection to mirror the saved selection state
$ws.Range("G95").Select() | Out-Null

$wb.Save() | Out-Null
